$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Replace "Working" (in "Week 12: Working on URS & Test plan.") with
#    "Agree on Major/Minor requirements. Work" and force it to live in
#    its own run (distinct from the "Week 12: " run before it and the
#    " on URS & Test plan." run after it) by nudging a character
#    property on the range - this causes the engine to persist the
#    range as a standalone <w:r>.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Working", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "Agree on Major/Minor requirements. Work"
$rng.Bold = 1
$rng.Bold = 0

# ---------------------------------------------------------------------
# 2) Drop a collapsed "_GoBack" bookmark right between "...Work" and
#    " on URS & Test plan.". Word only ever keeps a single "_GoBack"
#    bookmark, so adding it here automatically removes it from wherever
#    it used to be (the trailing empty paragraph at the end of the doc).
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("requirements. Work", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBack = $d.Range($rng2.End, $rng2.End)
$d.Bookmarks.Add("_GoBack", $goBack)
